$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of order data (Order_ID, Restaurant_Location, Customer_Location,
# Order_Time, Pickup_Time, Delivery_Time, Delivery_Distance_km, Driver_ID, Delivery_Status)
$data = @(
    @("ORD106", "South Street", "North Zone",   0.43402777777777773, 0.44444444444444442, 0.45833333333333331, 7.2, "D01", "Delivered"),
    @("ORD107", "Downtown",     "South Zone",   0.4375,               0.44791666666666669, 0.46180555555555558, 4.8, "D02", "Delivered"),
    @("ORD108", "East Zone",    "East Zone",    0.44097222222222227,  0.4513888888888889,  0.46875,              5,   "D03", "Delayed"),
    @("ORD109", "West End",     "West Zone",    0.44444444444444442,  0.4548611111111111,  0.47222222222222227,  6.5, "D04", "Delivered"),
    @("ORD110", "Central Mall", "Central Zone", 0.44791666666666669,  0.45833333333333331, 0.47569444444444442,  7.8, "D05", "Delivered"),
    @("ORD111", "North Park",   "North Zone",   0.4513888888888889,   0.46180555555555558, 0.47916666666666669,  5.0999999999999996, "D01", "Delayed"),
    @("ORD112", "South Street", "South Zone",   0.4548611111111111,   0.46527777777777773, 0.4861111111111111,   8.4, "D02", "Delivered"),
    @("ORD113", "Downtown",     "East Zone",    0.45833333333333331,  0.46875,              0.48958333333333331,  4.5, "D03", "Delivered"),
    @("ORD114", "East Zone",    "West Zone",    0.46180555555555558,  0.47222222222222227,  0.49652777777777773,  7.9, "D04", "Delayed"),
    @("ORD115", "West End",     "Central Zone", 0.46527777777777773,  0.47569444444444442,  0.5,                  6.2, "D05", "Delivered"),
    @("ORD116", "Central Mall", "North Zone",   0.46875,               0.47916666666666669, 0.50347222222222221,  5.8, "D01", "Delivered"),
    @("ORD117", "North Park",   "South Zone",   0.47222222222222227,  0.4826388888888889,  0.50694444444444442,  8,   "D02", "Delivered"),
    @("ORD118", "South Street", "East Zone",    0.47569444444444442,  0.4861111111111111,  0.51388888888888895,  7.3, "D03", "Delayed"),
    @("ORD119", "Downtown",     "West Zone",    0.47916666666666669,  0.48958333333333331, 0.51736111111111105,  4.5999999999999996, "D04", "Delivered"),
    @("ORD120", "East Zone",    "Central Zone", 0.4826388888888889,   0.49305555555555558, 0.52083333333333337,  6.9, "D05", "Delivered")
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).NumberFormat = "h:mm"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).NumberFormat = "h:mm"
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).NumberFormat = "h:mm"
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}

$ws.Range("A5").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G25").Select()
